# Updated Python Selenium Framework
# - Adds a new "BookingModule" worksheet (flight booking test data) after "LoginModule"
# - Tweaks a couple of existing LoginModule data cells
# - Updates the selected cell in both sheets

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LoginModule")

# ---------------------------------------------------------------------------
# 1. LoginModule tweaks
# ---------------------------------------------------------------------------
# Row 3 (negative-path test data): wrong username -> "mercury"
$ws1.Range("B3").Value = "mercury"
# Row 4 (negative-path test data): should not execute -> "No"
$ws1.Range("D4").Value = "No"

# Update the remembered selection on LoginModule
$ws1.Range("B7").Select()

# ---------------------------------------------------------------------------
# 2. Add the BookingModule worksheet right after LoginModule
# ---------------------------------------------------------------------------
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "BookingModule"
$ws2.Cells.Clear()

# ---------------------------------------------------------------------------
# 3. Formatting - reuse the existing header / body styles from LoginModule so
#    no redundant style entries are created.
# ---------------------------------------------------------------------------
# Header row (style used by LoginModule row 1)
$ws1.Range("A1").Copy()
$ws2.Range("A1:J1").PasteSpecial(-4122)

# Body rows - default body style (LoginModule row 2, columns A/B/D)
$ws1.Range("A2").Copy()
$ws2.Range("A2:J4").PasteSpecial(-4122)

# Body cells that use the alternate body style (LoginModule row 2, column C)
$ws1.Range("C2").Copy()
$ws2.Range("C2:C4").PasteSpecial(-4122)
$ws2.Range("E2:E4").PasteSpecial(-4122)
$ws2.Range("G3:H4").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 4. BookingModule data
# ---------------------------------------------------------------------------
# Header row
$ws2.Range("A1").Value = "TestCaseName"
$ws2.Range("B1").Value = "UserName"
$ws2.Range("C1").Value = "Password"
$ws2.Range("D1").Value = "No_of_Passengers"
$ws2.Range("E1").Value = "Departing_From"
$ws2.Range("F1").Value = "Departing_Day"
$ws2.Range("G1").Value = "Passenger_FirstName"
$ws2.Range("H1").Value = "Passenger_LastName"
$ws2.Range("I1").Value = "CreditCard_No"
$ws2.Range("J1").Value = "Execute"

# Row 2
$ws2.Range("A2").Value = "test_flightbooking"
$ws2.Range("B2").Value = "Pradyumna"
$ws2.Range("C2").Value = "mercury"
$ws2.Range("D2").Value = "1"
$ws2.Range("E2").Value = "Paris"
$ws2.Range("F2").Value = "18"
$ws2.Range("G2").Value = "Pradyumna"
$ws2.Range("H2").Value = "R"
$ws2.Range("I2").Value = "0123456789"
$ws2.Range("J2").Value = "Yes"

# Row 3
$ws2.Range("A3").Value = "test_flightbooking"
$ws2.Range("B3").Value = "mercury"
$ws2.Range("C3").Value = "mercury"
$ws2.Range("D3").Value = "1"
$ws2.Range("E3").Value = "Paris"
$ws2.Range("F3").Value = "18"
$ws2.Range("G3").Value = "Sandhya "
$ws2.Range("H3").Value = "C"
$ws2.Range("I3").Value = "0123456789"
$ws2.Range("J3").Value = "No"

# Row 4
$ws2.Range("A4").Value = "test_flightbooking"
$ws2.Range("B4").Value = "mercury"
$ws2.Range("C4").Value = "mercury"
$ws2.Range("D4").Value = "1"
$ws2.Range("E4").Value = "Paris"
$ws2.Range("F4").Value = "18"
$ws2.Range("G4").Value = "Samrudh"
$ws2.Range("H4").Value = "P"
$ws2.Range("I4").Value = "0123456789"
$ws2.Range("J4").Value = "No"

# ---------------------------------------------------------------------------
# 5. Column widths (best fit to content, mirrors the other sheet's layout)
# ---------------------------------------------------------------------------
for ($i = 1; $i -le 10; $i++) {
    $ws2.Columns.Item($i).AutoFit()
}

# ---------------------------------------------------------------------------
# 6. Page setup / selection / view state for the new sheet
# ---------------------------------------------------------------------------
$ws2.PageSetup.Orientation = 1
$ws2.Range("J6").Select()

# Keep LoginModule as the active/visible sheet, as before
$ws1.Activate()
